# NotlarIleriJavaGuz2017.xlsx - "Ileri java ders icerik" content update
#
# The grade sheet (rows 13-21 and 23-33) had blank "quiz" entries in
# columns T (Quiz1) and U (Quiz2) - these get explicit 0 scores entered,
# matching the rows above (3-12) and row 22 that already carry values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToFill = @(13, 14, 15, 16, 17, 18, 19, 20, 21, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33)

foreach ($r in $rowsToFill) {
    $ws.Cells.Item($r, 20).Value = 0   # column T
    $ws.Cells.Item($r, 21).Value = 0   # column U
}

# Leave the cursor where the author left it when saving.
$ws.Range("T28").Select()
